$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Set hours worked for the first task row; dependent formulas recalc automatically.
$ws.Range("H6").Value = 2

# Re-seat these merged header ranges so they relocate to the end of the
# mergeCells list (matches how Excel re-emits merges touched during editing).
$ws.Range("AZ4:BA4").UnMerge()
$ws.Range("AZ4:BA4").Merge()
$ws.Range("AO4:AP4").UnMerge()
$ws.Range("AO4:AP4").Merge()
$ws.Range("AR4:AS4").UnMerge()
$ws.Range("AR4:AS4").Merge()
$ws.Range("AU4:AV4").UnMerge()
$ws.Range("AU4:AV4").Merge()
$ws.Range("AX4:AY4").UnMerge()
$ws.Range("AX4:AY4").Merge()

# Update the active cell selection on the bottom-right frozen pane.
[void]$ws.Range("F6").Select()
